# icsd3-queries.xlsx — Sheet1 ("ICSD3-Queries") updates
#
# 1. Expand the Obstructive Sleep Apnea, Adult query with an extra OR clause.
# 2. Add a new "Chronic Insomnia Disorder" row (designation + query).
# 3. Row 2 no longer needs its taller custom height now that the sheet has
#    grown (let Excel drop back to the sheet's default row height).
# 4. Column B needs to widen to fit the longer text.
# 5. The active selection ends up on B10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Updated OSA query (adds "OR icsdadultosa_psggt15events" to the second clause)
$ws.Range("B2").Value = "((icdadultosa_symptoms OR  icdadultosa_history) AND (SUM(icsdadultosa_psg5events) >= 5)) OR ((SUM(icsdadultosa_psg15events) >= 15) OR icsdadultosa_psggt15events)"

# 2. New row for Chronic Insomnia Disorder
$ws.Range("A3").Value = "Chronic Insomnia Disorder"
$ws.Range("B3").Value = "icsdinsom_history AND icsdinsom_symptoms"

# 3. Row 2 reverts to the sheet's default height (no explicit ht="30" anymore)
$ws.Rows.Item(2).AutoFit()

# 4. Widen column B so the longer query text fits (~155.57 characters wide)
$ws.Columns.Item(2).ColumnWidth = 154.6

# 5. Move/record the selected cell as B10
$ws.Range("B10").Select()
